$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (column G) values computed for the 2024 save_data regen.
# (Strike# replaced by K; row 16 and row 50 already equalled 0 and are unchanged.)
$kValues = @{
    2 = 1
    3 = 0
    4 = 1
    5 = 1
    6 = 3
    7 = 3
    8 = 1
    9 = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 1
    17 = 0
    18 = 1
    19 = 2
    20 = 0
    21 = 0
    22 = 2
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 0
    28 = 2
    29 = 1
    30 = 1
    31 = 0
    32 = 1
    33 = 0
    34 = 1
    35 = 0
    36 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 0
    41 = 1
    42 = 2
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 2
    48 = 1
    49 = 1
    51 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Output "Updated $($kValues.Count) K (column G) cells"
